# Sulimn.xlsx - Alpha 2.2 edit
# Adds a "Rings" section (headers + 5 data rows) below the existing tables,
# tweaks a couple of column widths, and moves the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "Rings" section -------------------------------------------------

# Row 11: section title "Rings" (merged A11:B11), same look as the other
# section headers in row 1 (bold, 14pt, centered).
$ws.Range("A11:B11").Font.Bold = $true
$ws.Range("A11:B11").Font.Size = 14
$ws.Range("A11:B11").HorizontalAlignment = -4108
$ws.Range("A11:B11").VerticalAlignment = -4108
$ws.Range("A11").Value = "Rings"
$ws.Range("A11:B11").Merge()
$ws.Rows.Item(11).RowHeight = 18.75

# Row 12: column headers "Amount" / "Value" (bold, 12pt, centered) - same
# style as the other "Amount"/"Value" header row (row 2, columns M:N).
$ws.Range("A12:B12").Font.Bold = $true
$ws.Range("A12:B12").Font.Size = 12
$ws.Range("A12:B12").HorizontalAlignment = -4108
$ws.Range("A12:B12").VerticalAlignment = -4108
$ws.Range("A12").Value = "Amount"
$ws.Range("B12").Value = "Value"
$ws.Rows.Item(12).RowHeight = 15.75

# Rows 13-17: ring data. Column A = plain centered text, column B = number
# formatted with a thousands separator, also centered.
$ringData = @(
    @("1 of 1 Type", 2000),
    @("2 of 1 Type", 4000),
    @("1 of 2 Types", 4000),
    @("1 of 3 Types", 8000),
    @("1 of 4 Types", 16000)
)

$r = 13
foreach ($pair in $ringData) {
    $ws.Range("A$r").Value = $pair[0]
    $ws.Range("A$r").HorizontalAlignment = -4108
    $ws.Range("A$r").VerticalAlignment = -4108

    $ws.Range("B$r").Value = $pair[1]
    $ws.Range("B$r").NumberFormat = "#,##0"
    $ws.Range("B$r").HorizontalAlignment = -4108
    $ws.Range("B$r").VerticalAlignment = -4108

    $r++
}

# --- Column width tweak ----------------------------------------------------
# Columns A:B, D:E, G:H, J:K, M:N, P:Q shrink slightly.
foreach ($colRange in @("A:B", "D:E", "G:H", "J:K", "M:N", "P:Q")) {
    $ws.Range($colRange).ColumnWidth = 12.6
}

# --- Selection ---------------------------------------------------------
$ws.Range("S2").Select()
